# Add the new daily ranking sheet "2025-08-18" after the last existing sheet,
# populate it with the day's ranking data, and restore the original active sheet.

$wb = $excel.ActiveWorkbook
$ws57 = $wb.Worksheets.Item("2025-08-17")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-08-18"

# Match page margins used by the other daily ranking sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy the header row (values + bold/bordered style) from the previous day's sheet
$ws57.Range("A1:D1").Copy($newSheet.Range("A1:D1"))

$dataBlock = @"
1|願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜|ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)|第5話-2：魔法のお勉強
2|姫様“拷問”の時間です|原作:春原ロビンソン　漫画:ひらけい|拷問147
3|えろいことするために巨乳美少女奴隷を買ったはずが、お師匠さまと慕われて思った通りにいかなくなる話|佐藤36(作画) 煮豆シューター(原作)|第4話前半
4|宇崎ちゃんは遊びたい！|丈(著者)|第126話
5|窓際編集とバカにされた俺が、双子ＪＫと同居することになった|うさおとめ(著者) 茨木野(原作) トモゼロ(キャラクター原案)|第5話②
6|いとこのこ|いぬちく(著者)|第37話
7|みつばものがたり 呪いの少女と死の輪舞《ロンド》|堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)|第11話：勝利の美酒
8|序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する|作画：マエD 原作：新人|第5話(4)
9|りゅうとあまがみ|角丸柴朗(著者)|第二話・お肉は何処？①
10|悪役貴族として必要なそれ|まさこりん(原作) 夏野うみ(作画) 村カルキ(キャラクターデザイン)|第18話②
11|アイツノカノジョ|肉丸|第56話
12|最強の少年聖騎士、転生者を狩る|作画：御塩 原作：宇奈木ユラ|第7話(1)
13|最強勇者パーティーは愛が知りたい|山田肌襦袢|第29話「きみがきみであればいい」
14|おんなのこのけんをてにいれた|福岡太朗(著者)|18本目
15|ダメ人間の愛しかた|岩葉(著者)|第19話前編　ダメ人間と新生活の彼女
16|氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic|漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca|第9話
17|リビルドワールド|綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)|第72話②
18|サーシャちゃんとクラスメイトオタクくん|はぐはぐ(著者)|第83話
19|「才能の器」で目指す迷宮最深部 スキル横伸ばしのはずが、万能チートだった!|漫画：かくばやしつよし 原作：とんび キャラクター原案： りりんら|第40話
20|追放貴族は、外れスキル【古代召喚】で英霊たちと辺境領地を再興する～英霊たちを召喚したら慕われたので、最強領地を作り上げます～|三神みかみ(作画) たかた　ちひろ(原作)|第５話後半
21|魔都精兵のスレイブ|原作:タカヒロ　漫画:竹村洋平|第159話　神奴隷
22|異世界迷宮のオーパーツ|三狛ハル(著者)|第3話-②：おそらく高貴な布
23|聖液鍛冶屋のエロランタ|しげきっくす(著者)|第7話
24|半人前の恋人|川田大智|第50話
25|ハズレ職〈召喚士〉がS級万能職に化けました～無能と蔑まれた俺、伝説の召喚獣達に懐かれ力が覚醒したので世界最強です～|野呂まこと(作画) ヒツキノドカ(原作)|第4話前半
26|不徳のギルド|河添太一|第９７話：立派に育った所
27|よくわからないけれど異世界に転生していたようです|内々けやき あし カオミン|第137話 よくわからないけれど脱出するみたいです（２）
28|ぽんドロイド！ はまさん|はれやまはれぞう(著者)|第6話
29|リアリスト魔王による聖域なき異世界改革|鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)|第68幕②
30|世界の終わりの世界録(アンコール)|雨水龍(著者) 細音啓(原作) ふゆの春秋(キャラクター原案)|第96話②
31|今日から僕は、彼女の✕✕を解消する。|コアヤアコ(著者)|第1話
32|追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。|六志麻あさ 業務用餅 kisui|第７０話
33|アラサーがVTuberになった話。|犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)|第26話
34|貴族令嬢がジャンクフード食って「美味いですわ！」するだけの話|ごくげつ(作画) パイルバンカー串山(原作)|第4話前半
35|くらいあの子としたいこと|碇マナツ(著者)|第82話
36|俺堕ちスレイブヒーローコレクション|ゆっ栗栖(著者)|第12話前半
37|辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？|tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)|第10話
38|スキル【再生】と【破壊】から始まる最強冒険者ライフ～ごみ拾いと追放されたけど規格外の力で成り上がる！ ～|華尾ス太郎(作画) シンギョウガク(原作) Tea(キャラクター原案)|第4話前半
39|転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く|都尾琉(漫画) 御峰。(原作)|第27話③
40|小林さんちのメイドラゴン|クール教信者|第148話
41|ハーレムより平穏を！異世界で静かにニート姫させてくれ|さかたはるき(原作) かわやばぐ(作画)|第14話後半
42|治癒魔法の間違った使い方 ~戦場を駆ける回復要員~|九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)|第81話その3
43|きみの願いが叶うまで|浅月のりと(著者)|第4話-2
44|役目を果たした日陰の勇者は、辺境で自由に生きていきます|船野真帆(作画) 丘野優(原作) 布施龍太(キャラクター原案)|第5話後半
45|理想のヒモ生活|日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)|第86話　その3
46|おはらい箱の天才付与術師は、辺境で悠々自適に暮らしたい～万能付与術で気付いたら辺境が世界最強の快適拠点になっていた～|柊木 楸(作画) 水無月(原作) 布施龍太(キャラクター原案)|第4話前半
47|愚かな天使は悪魔と踊る|アズマサワヨシ(著者)|第101話②
48|追放された元雑用係、規格外の技術で「最高の修繕師」と呼ばれるようになりました～SSSランクパーティーや王族からの依頼が止まりません～|日高(作画) あざね(原作)|第４話前半
49|うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。|三上こた こばやし少女 寝子空兄 ゆがー|第1話　災禍の悪夢
50|クラスで２番目に可愛い女の子と友だちになった|尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)|第34話②
"@

$lines = $dataBlock -split "`n"
$rowIndex = 2
foreach ($line in $lines) {
    if ($line.Trim().Length -eq 0) { continue }
    $parts = $line -split "\|"
    $newSheet.Cells.Item($rowIndex, 1).Value = [int]$parts[0]
    $newSheet.Cells.Item($rowIndex, 2).Value = $parts[1]
    $newSheet.Cells.Item($rowIndex, 3).Value = $parts[2]
    $newSheet.Cells.Item($rowIndex, 4).Value = $parts[3]
    $rowIndex++
}

# Restore Sheet1 as the active sheet (matches the workbook's original active-tab state)
$wb.Worksheets.Item(1).Activate()
